# Applies the edit described by the diff:
#  1. Updates the "Noisy"/"NLM-LBP" numeric values on the sigma_010 and
#     sigma_025 sheets with their newly recomputed (higher precision) values.
#  2. Adds a new worksheet "sigma_050" (with the same 3-column layout:
#     Rows / Noisy / NLM-LBP, 10 data rows + a "Média" summary row) at the
#     end of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update existing sheets sigma_010 (sheet2) and sigma_025 (sheet3)
# ---------------------------------------------------------------------------

$sigma010 = $wb.Worksheets.Item("sigma_010")
$sigma025 = $wb.Worksheets.Item("sigma_025")

# Row => (B value, C value) for sigma_010
$sigma010Updates = @{
    2  = @(28.20536939091576, 33.08503803106901)
    3  = @(28.16086738384023, 33.08837833079527)
    4  = @(28.20350085034629, 33.06534372179797)
    5  = @(28.16969945537825, 33.08270338276674)
    6  = @(28.17414451097907, 33.10213521319915)
    7  = @(28.18914747498484, 33.09174468346005)
    8  = @(28.21107258196324, 33.09414581669628)
    9  = @(28.19619817954242, 33.10478961634517)
    10 = @(28.16780142480074, 33.10719441765799)
    11 = @(28.18455121027186, 33.07197546099432)
    12 = @(28.18623524630227, 33.0893448674782)
}

foreach ($row in $sigma010Updates.Keys) {
    $vals = $sigma010Updates[$row]
    $sigma010.Cells.Item($row, 2).Value = $vals[0]
    $sigma010.Cells.Item($row, 3).Value = $vals[1]
}

# Row => (B value, C value) for sigma_025
$sigma025Updates = @{
    2  = @(19.68219233886653, 30.10562828032243)
    3  = @(19.6540882895902,  30.07341050112868)
    4  = @(19.66613858350984, 30.08259349767796)
    5  = @(19.6677116555131,  30.09300099501445)
    6  = @(19.64267509982664, 30.07468897153784)
    7  = @(19.67092763213526, 30.11529876680587)
    8  = @(19.67087660083289, 30.09790456347201)
    9  = @(19.65637676892247, 30.14298275179414)
    10 = @(19.66005187916833, 30.08981100428429)
    11 = @(19.66410222091681, 30.13501810845969)
    12 = @(19.66351410692821, 30.10103374404974)
}

foreach ($row in $sigma025Updates.Keys) {
    $vals = $sigma025Updates[$row]
    $sigma025.Cells.Item($row, 2).Value = $vals[0]
    $sigma025.Cells.Item($row, 3).Value = $vals[1]
}

# ---------------------------------------------------------------------------
# 2. Add new sheet "sigma_050" at the end of the workbook
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sigma050 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sigma050.Name = "sigma_050"

# Header row
$sigma050.Cells.Item(1, 1).Value = "Rows"
$sigma050.Cells.Item(1, 2).Value = "Noisy"
$sigma050.Cells.Item(1, 3).Value = "NLM-LBP"

# Data rows: Row index => (A, B, C)
$sigma050Data = @{
    2  = @(0, 14.48338446106393, 25.17579886182194)
    3  = @(1, 14.48533510203679, 25.20598679617718)
    4  = @(2, 14.48949615537303, 25.21476546815518)
    5  = @(3, 14.48853078166531, 25.23304109178026)
    6  = @(4, 14.49530774530084, 25.28531228590264)
    7  = @(5, 14.47969177291145, 25.23009284335531)
    8  = @(6, 14.50279262769881, 25.23244070684862)
    9  = @(7, 14.48597213448632, 25.1715663479612)
    10 = @(8, 14.4856828025496,  25.29536754363697)
    11 = @(9, 14.48139923490796, 25.24020234646948)
}

foreach ($row in $sigma050Data.Keys) {
    $vals = $sigma050Data[$row]
    $sigma050.Cells.Item($row, 1).Value = $vals[0]
    $sigma050.Cells.Item($row, 2).Value = $vals[1]
    $sigma050.Cells.Item($row, 3).Value = $vals[2]
}

# "Média" summary row
$sigma050.Cells.Item(12, 1).Value = "Média"
$sigma050.Cells.Item(12, 2).Value = 14.4877592817994
$sigma050.Cells.Item(12, 3).Value = 25.22845742921088

Write-Host "Edit applied: sigma_010/sigma_025 values refreshed, sigma_050 sheet added."
